# Edit: "Fuel" sheet update for both CH and SIN database
# - updated PEN & CO2 in "Fuel" sheet for various fuels CH (data source: KBOB 2009/1:2016)
# - updated PEN & CO2 for natural gas (NG) in "Fuel" sheet for SIN (data source: ecoinvent 3.4 ...)

$wb = $excel.ActiveWorkbook

$wsDHW   = $wb.Worksheets.Item("DHW")
$wsHEAT  = $wb.Worksheets.Item("HEATING")
$wsCOOL  = $wb.Worksheets.Item("COOLING")
$wsELEC  = $wb.Worksheets.Item("ELECTRICITY")
$wsFUEL  = $wb.Worksheets.Item("FUELS")

# ---------------------------------------------------------------------------
# ELECTRICITY sheet ("Swiss consumer energy mix" / GRID row) - updated PEN/CO2
# ---------------------------------------------------------------------------
$wsELEC.Range("E3").Value = 2.52
$wsELEC.Range("F3").Value = 0.028
$wsELEC.Range("H3").Value = "KBOB 2009/1:2016, ID 45.020 CH-Verbrauchermix, costs in USD-2015"

# ---------------------------------------------------------------------------
# FUELS sheet - updated PEN/CO2 + references for Natural Gas, Electricity,
# Oil, Coal, Wood, and a brand-new "Biogas" row
# ---------------------------------------------------------------------------

# Natural Gas (row 2)
$wsFUEL.Range("C2").Value = 1.06
$wsFUEL.Range("D2").Formula = "=0.228/3.6"
$wsFUEL.Range("F2").Value = "KBOB 2009/1:2016, ID 41.002 Erdgas"

# Electricity (row 3)
$wsFUEL.Range("C3").Value = 2.52
$wsFUEL.Range("D3").Formula = "=0.102/3.6"
$wsFUEL.Range("F3").Value = "KBOB 2009/1:2016, ID 45.020 CH-Verbrauchermix"

# Oil (row 5)
$wsFUEL.Range("C5").Value = 1.23
$wsFUEL.Range("D5").Formula = "=0.301/3.6"
$wsFUEL.Range("F5").Value = "KBOB 2009/1:2016, ID 41.001 Heizöl"

# Coal (row 6)
$wsFUEL.Range("C6").Formula = "=(1.2+1.45)/2"
$wsFUEL.Range("D6").Formula = "=((0.399+0.439)/2)/3.6"
$wsFUEL.Range("F6").Value = "KBOB 2009/1:2016, ID 41.004/41.005 (average)"

# Wood (row 7)
$wsFUEL.Range("C7").Value = 0.11600000000000001
$wsFUEL.Range("D7").Formula = "=0.027/3.6"
$wsFUEL.Range("F7").Value = "KBOB 2009/1:2016, ID 41.006 Stückholz (average)"

# New Biogas row (row 8)
$wsFUEL.Range("A7:F7").Copy() | Out-Null
$wsFUEL.Range("A8:F8").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$wsFUEL.Range("A8").Value = "Biogas"
$wsFUEL.Range("B8").Value = "BIOGAS"
$wsFUEL.Range("C8").Value = 0.29899999999999999
$wsFUEL.Range("D8").Formula = "=0.13/3.6"
$wsFUEL.Range("E8").Value = ""
$wsFUEL.Range("F8").Value = "KBOB 2009/1:2016, ID 41.009 Biogas"

$wsFUEL.Range("A9:F9").Value = ""

$wsFUEL.PageSetup.PaperSize = 9
$wsFUEL.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Restore the selections that were active on each sheet, finishing with the
# FUELS sheet active (as it was when the workbook was last saved).
# ---------------------------------------------------------------------------
$wsDHW.Activate()
$wsDHW.Range("A7").Select() | Out-Null

$wsHEAT.Activate()
$wsHEAT.Range("A7").Select() | Out-Null

$wsELEC.Activate()
$wsELEC.Range("E4").Select() | Out-Null

$wsFUEL.Activate()
$excel.ActiveWindow.Zoom = 150
$wsFUEL.Range("C16").Select() | Out-Null
